# Creating common utility for RunModes
#
# Swap the "y"/"n" RunMode flags between the two TestCases rows
# (AddCustomerTest / OpenAccountTest), then leave the workbook positioned
# the way it was when last saved: cursor on D7 of TestCases, and the
# TestData sheet active/selected with its cursor on C10.

$wb  = $excel.ActiveWorkbook
$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestData  = $wb.Worksheets.Item("TestData")

# B2 (AddCustomerTest RunMode) was "y", becomes "n"
# B3 (OpenAccountTest RunMode) was "n", becomes "y"
$wsTestCases.Range("B2").Value = "n"
$wsTestCases.Range("B3").Value = "y"

# Update the remembered selection on TestCases ...
$wsTestCases.Range("D7").Select()

# ... then switch to / activate the TestData sheet, keeping its own
# selection where it already was (C10).
$wsTestData.Activate()
$wsTestData.Range("C10").Select()
